$d = $word.ActiveDocument

# The document currently ends with the "Hola" paragraph, immediately
# followed by the sectPr. We need to append, after "Hola":
#   (empty)
#   "Capitulo 1"   (with spellcheck proofErr markers around "Capitulo")
#   (empty)
#   "Título 1"
#   (empty)
#   "Tema 1"
# All new paragraphs keep the same es-MX paragraph-mark run formatting
# used throughout the rest of the document.

# First create one fresh empty paragraph at the end of the story so we
# have a collapsed insertion point to target with InsertXML (this avoids
# disturbing the existing "Hola" paragraph's own run).
$count = $d.Paragraphs.Count
$last = $d.Paragraphs($count)
$last.Range.InsertParagraphAfter()

# Re-fetch the newly created (currently empty) last paragraph.
$count = $d.Paragraphs.Count
$last = $d.Paragraphs($count)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$langRpr = '<w:rPr><w:lang w:val="es-MX"/></w:rPr>'
$pPr = "<w:pPr>$langRpr</w:pPr>"

$pEmpty = "<w:p $wNs>$pPr</w:p>"
$pCapitulo = "<w:p $wNs>$pPr" + `
    '<w:proofErr w:type="spellStart"/>' + `
    "<w:r>$langRpr<w:t>Capitulo</w:t></w:r>" + `
    '<w:proofErr w:type="spellEnd"/>' + `
    "<w:r>$langRpr<w:t xml:space=`"preserve`"> 1</w:t></w:r>" + `
    '</w:p>'
$pTitulo = "<w:p $wNs>$pPr<w:r>$langRpr<w:t>Título 1</w:t></w:r></w:p>"
$pTema = "<w:p $wNs>$pPr<w:r>$langRpr<w:t>Tema 1</w:t></w:r></w:p>"

# InsertXML on a collapsed range at an (empty) paragraph fills that
# paragraph with the first fragment paragraph's content and appends the
# remaining fragment paragraphs as new paragraphs right after it - giving
# us exactly the six target paragraphs in one shot.
$xml = $pEmpty + $pCapitulo + $pEmpty + $pTitulo + $pEmpty + $pTema
$last.Range.InsertXML($xml)

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
